$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency price/volume data (scraped update) and apply two
# row re-orderings: rows 34/35 (ARBITRUM <-> ImmutableX) and rows 47/48 (EnergySwap <-> PaxDollar).

# Force plain-numeric-looking Price values to stay as text, matching the
# original inline-string cell type (avoids Excel auto-converting to numbers).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.093.45"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "1.887.35"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "307.42"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.5156"
$ws.Range("D8").Value = "0.3723"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("D9").Value = "0.07211"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").Value = "0.9024"
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("D11").Value = "21.02"
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("D12").Value = "0.07624"
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("D13").Value = "1.893.65"
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("D14").Value = "94.53"
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").Value = "5.266"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("E18").Value = "  +2.34%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "27.145.90"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").Value = "5.056"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").Value = "2.136.68"
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("D23").Value = "10.56"
$ws.Range("E23").Value = "  +2.35%  "
$ws.Range("D24").Value = "6.421"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "146.45"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "18.02"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("D28").Value = "2.171"
$ws.Range("E28").Value = "  +5.57%  "
$ws.Range("D29").Value = "114.59"
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("D30").Value = "4.980"
$ws.Range("E30").Value = "  +6.85%  "
$ws.Range("D31").Value = "4.812"
$ws.Range("E31").Value = "  +3.91%  "
$ws.Range("D32").Value = "0.09205"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").Value = "0.05066"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.7626"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.196"
$ws.Range("E35").Value = "  +4.57%  "
$ws.Range("D36").Value = "2.978"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "3.273"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("D38").Value = "2.580"
$ws.Range("E38").Value = "  +3.07%  "
$ws.Range("D39").Value = "0.5625"
$ws.Range("E39").Value = "  +5.66%  "
$ws.Range("D40").Value = "0.01993"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").Value = "1.075"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").Value = "8.994"
$ws.Range("E42").Value = "  +7.42%  "
$ws.Range("D43").Value = "118.54"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "6.576"
$ws.Range("E44").Value = "  +1.77%  "
$ws.Range("D45").Value = "0.1506"
$ws.Range("E45").Value = "  +3.39%  "
$ws.Range("D46").Value = "0.4801"
$ws.Range("E46").Value = "  +3.42%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "0.9996"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "10.15"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("D49").Value = "1.586"
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("D50").Value = "37.20"
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("D51").Value = "63.72"
$ws.Range("E51").Value = "  +1.40%  "
